# Insert a new data row at row 690 (pushing existing rows 690-788 down to 691-789)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(690).Insert()

$ws.Cells.Item(690, 1).Value2  = 8
$ws.Cells.Item(690, 2).Value2  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(690, 3).Value2  = 'Coquimbo'
$ws.Cells.Item(690, 4).Value2  = 44984
$ws.Cells.Item(690, 5).Value2  = 4
$ws.Cells.Item(690, 6).Value2  = 100112024
$ws.Cells.Item(690, 7).Value2  = 'Choclo'
$ws.Cells.Item(690, 8).Value2  = 'Choclero'
$ws.Cells.Item(690, 9).Value2  = 'Primera'
$ws.Cells.Item(690, 10).Value2 = 8000
$ws.Cells.Item(690, 11).Value2 = 450
$ws.Cells.Item(690, 12).Value2 = 500
$ws.Cells.Item(690, 13).Value2 = 475
$ws.Cells.Item(690, 14).Value2 = '$/unidad'
$ws.Cells.Item(690, 15).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(690, 16).Value2 = 475
$ws.Cells.Item(690, 17).Value2 = 1
$ws.Cells.Item(690, 18).Value2 = 'Hortaliza'
